# fix: sexting intensity escalation + PPV0 phase detection bug across all 23 models
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EmilyBellJourney")

$ws.Range("B4").Value = "cum with me babe... right now, watch 🥵"
$ws.Range("B6").Value = "I'm cumming right now... don't miss this"
$ws.Range("B7").Value = "I'm SO close babe... wait for me, I want you to watch the second it happens 😏"
$ws.Range("B8").Value = "oh fuck oh fuck 🥵"
$ws.Range("B9").Value = "you're not ready for this one babe... but I need you to see it 🥵"
$ws.Range("B11").Value = "I'm about to lose it and I need you to see what's happening to me right now"
$ws.Range("B12").Value = "my fingers are deep inside me and I can't stop moaning... I hope my neighbors can't hear this 😏"
$ws.Range("B13").Value = "I'm rubbing my clit so fast right now and god it feels so good thinking about you watching"
$ws.Range("B14").Value = "fuckkk 🥵"
$ws.Range("B15").Value = "look at me... this is ALL because of you and I can't stop 🥵"
$ws.Range("B17").Value = "what do you want me to do next? seriously I'll do literally anything you tell me right now"
$ws.Range("B18").Value = "I need your hands on every part of me right now babe... I keep imagining it and my body is going crazy 🥵"
$ws.Range("B19").Value = "but I can't stop now... my fingers slipped inside and I'm soaking wet because of you"
$ws.Range("B20").Value = "oh wow... okay I did NOT expect to feel like this 😏"
$ws.Range("B21").Value = "oh my god I can't believe I'm sending this... but you need to see what you did 🥵"
$ws.Range("B23").Value = "okay I'm definitely touching myself right now and I blame you entirely babe 🥵"
$ws.Range("B24").Value = "I wasn't planning on going there tonight but you're literally making me so wet I can't think straight"
$ws.Range("B25").Value = "soo you liked that huh? because honestly my heart is racing knowing you just saw that 😏"
